# Update "想去人数" (interest count) figures on the 展览, 演出, and 全部类型 sheets
# to match the freshly generated data snapshot (gh-pages output).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 25
$ws.Range("F3").Value = 482
$ws.Range("F4").Value = 217
$ws.Range("F5").Value = 83
$ws.Range("F6").Value = 325
$ws.Range("F8").Value = 518
$ws.Range("F9").Value = 108
$ws.Range("F10").Value = 1330
$ws.Range("F12").Value = 1096
$ws.Range("F13").Value = 172
$ws.Range("F15").Value = 292
$ws.Range("F16").Value = 118
$ws.Range("F17").Value = 255
$ws.Range("F18").Value = 1681
$ws.Range("F20").Value = 272
$ws.Range("F21").Value = 252
$ws.Range("F22").Value = 2821
$ws.Range("F27").Value = 1218
$ws.Range("F29").Value = 2848
$ws.Range("F30").Value = 1646
$ws.Range("F31").Value = 87
$ws.Range("F32").Value = 124
$ws.Range("F33").Value = 685
$ws.Range("F35").Value = 1899
$ws.Range("F36").Value = 902
$ws.Range("F37").Value = 1907
$ws.Range("F39").Value = 32
$ws.Range("F40").Value = 34
$ws.Range("F42").Value = 49
$ws.Range("F43").Value = 898
$ws.Range("F44").Value = 810
$ws.Range("F45").Value = 1045
$ws.Range("F46").Value = 125
$ws.Range("F47").Value = 448
$ws.Range("F48").Value = 227
$ws.Range("F49").Value = 3366

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 12
$ws.Range("F12").Value = 808
$ws.Range("F18").Value = 9

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 482
$ws.Range("F3").Value = 217
$ws.Range("F4").Value = 83
$ws.Range("F5").Value = 12
$ws.Range("F7").Value = 325
$ws.Range("F9").Value = 518
$ws.Range("F10").Value = 108
$ws.Range("F11").Value = 1330
$ws.Range("F13").Value = 1096
$ws.Range("F14").Value = 172
$ws.Range("F16").Value = 292
$ws.Range("F17").Value = 118
$ws.Range("F18").Value = 255
$ws.Range("F19").Value = 1681
$ws.Range("F21").Value = 272
$ws.Range("F22").Value = 252
$ws.Range("F23").Value = 2822
$ws.Range("F27").Value = 1218
$ws.Range("F28").Value = 2848
$ws.Range("F29").Value = 1646
$ws.Range("F30").Value = 87
$ws.Range("F32").Value = 124
$ws.Range("F33").Value = 808
$ws.Range("F36").Value = 1899
$ws.Range("F38").Value = 902
$ws.Range("F39").Value = 1907
$ws.Range("F41").Value = 898
$ws.Range("F42").Value = 810
$ws.Range("F43").Value = 1045
$ws.Range("F44").Value = 125
$ws.Range("F45").Value = 448
$ws.Range("F47").Value = 227
$ws.Range("F48").Value = 3366

